$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FS")
$ws.Range("C6").Value = 0.4
$ws.Range("D6").Value = 0.4
$ws.Range("G6").Value = 36.84210526315789
$ws.Range("N6").Value = 1.1695906432748537
$ws.Range("D7").Value = 0.875
$ws.Range("E7").Value = 0.125
$ws.Range("D8").Value = 0.5
$ws.Range("E8").Value = 0.5
$ws.Range("K8").Value = 0.037037037037037035
$ws.Range("L8").Value = 0.9629629629629629
$ws.Range("C12").Value = 0.8
$ws.Range("E12").Value = 0.2
$ws.Range("G12").Value = 15.789473684210526
$ws.Range("N12").Value = 0.5847953216374269
$ws.Range("C13").Value = 0.0
$ws.Range("D13").Value = 1.0
$ws.Range("E13").Value = 0.0
$ws.Range("C14").Value = 0.3333333333333333
$ws.Range("E14").Value = 0.6666666666666666
$ws.Range("K14").Value = 0.018518518518518517
$ws.Range("L14").Value = 0.9814814814814815
$ws.Range("C18").Value = 0.4
$ws.Range("D18").Value = 0.6
$ws.Range("G18").Value = 31.57894736842105
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.022222222222222223
$ws.Range("N18").Value = 1.1695906432748537
$ws.Range("C20").Value = 0.16666666666666666
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.5
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.8
$ws.Range("K24").Value = 0.15555555555555556
$ws.Range("L24").Value = 0.044444444444444446
$ws.Range("N24").Value = 16.95906432748538
$ws.Range("D25").Value = 0.875
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.05555555555555555
$ws.Range("K25").Value = 0.8611111111111112
$ws.Range("L25").Value = 0.08333333333333333
$ws.Range("C26").Value = 0.8333333333333334
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.0
$ws.Range("J26").Value = 0.05555555555555555
$ws.Range("K26").Value = 0.12962962962962962
$ws.Range("L26").Value = 0.8148148148148148

$ws = $wb.Worksheets.Item("IF")
$ws.Range("C6").Value = 1.0
$ws.Range("D6").Value = 0.0
$ws.Range("G6").Value = 15.789473684210526
$ws.Range("J6").Value = 0.8666666666666667
$ws.Range("K6").Value = 0.08888888888888889
$ws.Range("N6").Value = 12.280701754385964
$ws.Range("C7").Value = 0.0
$ws.Range("E7").Value = 0.125
$ws.Range("J7").Value = 0.05555555555555555
$ws.Range("K7").Value = 0.875
$ws.Range("L7").Value = 0.06944444444444445
$ws.Range("C8").Value = 0.16666666666666666
$ws.Range("D8").Value = 0.16666666666666666
$ws.Range("E8").Value = 0.6666666666666666
$ws.Range("J8").Value = 0.018518518518518517
$ws.Range("K8").Value = 0.09259259259259259
$ws.Range("L8").Value = 0.8888888888888888
$ws.Range("C12").Value = 0.4
$ws.Range("D12").Value = 0.6
$ws.Range("G12").Value = 47.368421052631575
$ws.Range("J12").Value = 0.9555555555555556
$ws.Range("L12").Value = 0.044444444444444446
$ws.Range("N12").Value = 9.941520467836257
$ws.Range("J13").Value = 0.027777777777777776
$ws.Range("L13").Value = 0.041666666666666664
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.5
$ws.Range("J14").Value = 0.05555555555555555
$ws.Range("K14").Value = 0.12962962962962962
$ws.Range("L14").Value = 0.8148148148148148
$ws.Range("C18").Value = 0.6
$ws.Range("D18").Value = 0.4
$ws.Range("G18").Value = 57.89473684210527
$ws.Range("J18").Value = 0.7777777777777778
$ws.Range("K18").Value = 0.2222222222222222
$ws.Range("N18").Value = 18.71345029239766
$ws.Range("C19").Value = 0.75
$ws.Range("D19").Value = 0.25
$ws.Range("J19").Value = 0.1527777777777778
$ws.Range("K19").Value = 0.7916666666666666
$ws.Range("L19").Value = 0.05555555555555555
$ws.Range("C20").Value = 0.3333333333333333
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("J20").Value = 0.037037037037037035
$ws.Range("K20").Value = 0.09259259259259259
$ws.Range("L20").Value = 0.8703703703703703
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.4
$ws.Range("G24").Value = 47.368421052631575
$ws.Range("J24").Value = 0.7333333333333333
$ws.Range("K24").Value = 0.17777777777777778
$ws.Range("L24").Value = 0.08888888888888889
$ws.Range("N24").Value = 19.298245614035086
$ws.Range("D25").Value = 0.375
$ws.Range("E25").Value = 0.625
$ws.Range("J25").Value = 0.041666666666666664
$ws.Range("L25").Value = 0.09722222222222222
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.16666666666666666
$ws.Range("L26").Value = 0.7962962962962963

$ws = $wb.Worksheets.Item("IA")
$ws.Range("C6").Value = 0.8
$ws.Range("D6").Value = 0.0
$ws.Range("G6").Value = 26.31578947368421
$ws.Range("N6").Value = 1.7543859649122806
$ws.Range("C7").Value = 0.125
$ws.Range("D7").Value = 0.5
$ws.Range("J7").Value = 0.0
$ws.Range("K7").Value = 1.0
$ws.Range("D8").Value = 0.0
$ws.Range("E8").Value = 1.0
$ws.Range("C12").Value = 0.8
$ws.Range("E12").Value = 0.2
$ws.Range("G12").Value = 36.84210526315789
$ws.Range("N12").Value = 0.5847953216374269
$ws.Range("C14").Value = 0.16666666666666666
$ws.Range("D14").Value = 0.3333333333333333
$ws.Range("K14").Value = 0.018518518518518517
$ws.Range("L14").Value = 0.9814814814814815
$ws.Range("D18").Value = 0.4
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 26.31578947368421
$ws.Range("D20").Value = 0.5
$ws.Range("E20").Value = 0.5
$ws.Range("C24").Value = 0.8
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("N24").Value = 4.093567251461988
$ws.Range("C25").Value = 0.0
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.013888888888888888
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.013888888888888888
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666

$ws = $wb.Worksheets.Item("FS-IF")
$ws.Range("C6").Value = 0.2
$ws.Range("D6").Value = 0.0
$ws.Range("E6").Value = 0.8
$ws.Range("N6").Value = 1.1695906432748537
$ws.Range("C7").Value = 0.25
$ws.Range("D7").Value = 0.375
$ws.Range("E7").Value = 0.375
$ws.Range("K7").Value = 0.9861111111111112
$ws.Range("L7").Value = 0.013888888888888888
$ws.Range("C8").Value = 0.0
$ws.Range("D8").Value = 0.16666666666666666
$ws.Range("E8").Value = 0.8333333333333334
$ws.Range("C12").Value = 0.6
$ws.Range("D12").Value = 0.4
$ws.Range("E12").Value = 0.0
$ws.Range("D14").Value = 0.6666666666666666
$ws.Range("E14").Value = 0.3333333333333333
$ws.Range("D18").Value = 0.6
$ws.Range("E18").Value = 0.0
$ws.Range("J18").Value = 0.9555555555555556
$ws.Range("K18").Value = 0.044444444444444446
$ws.Range("N18").Value = 3.508771929824561
$ws.Range("C19").Value = 0.0
$ws.Range("D19").Value = 0.75
$ws.Range("J19").Value = 0.05555555555555555
$ws.Range("K19").Value = 0.9444444444444444
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("C24").Value = 0.6
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 63.1578947368421
$ws.Range("J24").Value = 0.8888888888888888
$ws.Range("K24").Value = 0.08888888888888889
$ws.Range("N24").Value = 14.619883040935672
$ws.Range("C25").Value = 0.625
$ws.Range("D25").Value = 0.375
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.013888888888888888
$ws.Range("K25").Value = 0.9305555555555556
$ws.Range("L25").Value = 0.05555555555555555
$ws.Range("C26").Value = 0.5
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.16666666666666666
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.24074074074074073
$ws.Range("L26").Value = 0.7222222222222222

$ws = $wb.Worksheets.Item("FS-IA")
$ws.Range("C6").Value = 1.0
$ws.Range("D6").Value = 0.0
$ws.Range("E6").Value = 0.0
$ws.Range("G6").Value = 10.526315789473683
$ws.Range("J6").Value = 1.0
$ws.Range("K6").Value = 0.0
$ws.Range("N6").Value = 1.1695906432748537
$ws.Range("C7").Value = 0.0
$ws.Range("D7").Value = 0.875
$ws.Range("J7").Value = 0.0
$ws.Range("K7").Value = 1.0
$ws.Range("D12").Value = 0.0
$ws.Range("E12").Value = 0.2
$ws.Range("G12").Value = 36.84210526315789
$ws.Range("C13").Value = 0.0
$ws.Range("E13").Value = 0.375
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.5
$ws.Range("C18").Value = 1.0
$ws.Range("D18").Value = 0.0
$ws.Range("G18").Value = 26.31578947368421
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.022222222222222223
$ws.Range("N18").Value = 1.1695906432748537
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.013888888888888888
$ws.Range("K19").Value = 0.9861111111111112
$ws.Range("D20").Value = 0.5
$ws.Range("E20").Value = 0.5
$ws.Range("K20").Value = 0.0
$ws.Range("L20").Value = 1.0
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.2
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 2.3391812865497075
$ws.Range("K25").Value = 1.0
$ws.Range("L25").Value = 0.0
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.018518518518518517
$ws.Range("L26").Value = 0.9444444444444444

$ws = $wb.Worksheets.Item("IF-IA")
$ws.Range("C6").Value = 0.6
$ws.Range("D6").Value = 0.4
$ws.Range("G6").Value = 36.84210526315789
$ws.Range("N6").Value = 0.0
$ws.Range("D7").Value = 0.875
$ws.Range("E7").Value = 0.125
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.0
$ws.Range("C8").Value = 0.0
$ws.Range("D8").Value = 0.6666666666666666
$ws.Range("K8").Value = 0.0
$ws.Range("L8").Value = 1.0
$ws.Range("C12").Value = 1.0
$ws.Range("E12").Value = 0.0
$ws.Range("G12").Value = 26.31578947368421
$ws.Range("C13").Value = 0.125
$ws.Range("D13").Value = 0.5
$ws.Range("E13").Value = 0.375
$ws.Range("C14").Value = 0.16666666666666666
$ws.Range("D14").Value = 0.0
$ws.Range("C18").Value = 0.4
$ws.Range("D18").Value = 0.4
$ws.Range("E18").Value = 0.2
$ws.Range("G18").Value = 52.63157894736842
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.022222222222222223
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 0.125
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.0
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.8333333333333334
$ws.Range("E20").Value = 0.16666666666666666
$ws.Range("C24").Value = 0.2
$ws.Range("E24").Value = 0.4
$ws.Range("G24").Value = 63.1578947368421
$ws.Range("J24").Value = 0.9111111111111111
$ws.Range("K24").Value = 0.022222222222222223
$ws.Range("L24").Value = 0.06666666666666667
$ws.Range("N24").Value = 10.526315789473683
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.25
$ws.Range("E25").Value = 0.625
$ws.Range("K25").Value = 0.9305555555555556
$ws.Range("L25").Value = 0.06944444444444445
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.05555555555555555
$ws.Range("K26").Value = 0.1111111111111111
$ws.Range("L26").Value = 0.8333333333333334

$ws = $wb.Worksheets.Item("FS-IF-IA")
$ws.Range("C6").Value = 0.0
$ws.Range("D6").Value = 0.6
$ws.Range("G6").Value = 63.1578947368421
$ws.Range("N6").Value = 1.7543859649122806
$ws.Range("C7").Value = 0.125
$ws.Range("D7").Value = 0.375
$ws.Range("E7").Value = 0.5
$ws.Range("K7").Value = 0.9861111111111112
$ws.Range("L7").Value = 0.013888888888888888
$ws.Range("C8").Value = 0.16666666666666666
$ws.Range("D8").Value = 0.16666666666666666
$ws.Range("K8").Value = 0.037037037037037035
$ws.Range("L8").Value = 0.9629629629629629
$ws.Range("C12").Value = 0.6
$ws.Range("E12").Value = 0.4
$ws.Range("C13").Value = 0.0
$ws.Range("D13").Value = 0.625
$ws.Range("E13").Value = 0.375
$ws.Range("C18").Value = 0.2
$ws.Range("D18").Value = 0.6
$ws.Range("E18").Value = 0.2
$ws.Range("G18").Value = 42.10526315789473
$ws.Range("J18").Value = 0.9777777777777777
$ws.Range("K18").Value = 0.022222222222222223
$ws.Range("N18").Value = 0.5847953216374269
$ws.Range("C19").Value = 0.0
$ws.Range("D19").Value = 1.0
$ws.Range("E19").Value = 0.0
$ws.Range("D20").Value = 0.6666666666666666
$ws.Range("E20").Value = 0.3333333333333333
$ws.Range("C24").Value = 0.8
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.9555555555555556
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.044444444444444446
$ws.Range("N24").Value = 5.263157894736842
$ws.Range("C25").Value = 0.0
$ws.Range("E25").Value = 0.25
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.6666666666666666
$ws.Range("E26").Value = 0.3333333333333333
$ws.Range("K26").Value = 0.07407407407407407
$ws.Range("L26").Value = 0.9074074074074074
